# "recreated docs w/o LINCS"
# Remove the "LINCS" Heading5 sub-section (its heading paragraph plus the
# "This is Agency common control..." body paragraph that immediately follows
# it) from the MP-1 control's common-control breakdown. Everything after it
# (MP-2, MP-6, MP-7, and their AWS sub-sections) shifts up to fill the gap;
# bookmark ids renumber automatically on save.

$d = $word.ActiveDocument

# Locate the "LINCS" heading paragraph via its bookmark so the edit is
# resilient to any incidental paragraph-index drift elsewhere in the doc.
$bm = $d.Bookmarks("lincs")
$headingPara = $bm.Range.Paragraphs(1)
$bodyPara = $headingPara.Next()

$start = $headingPara.Range.Start
$end = $bodyPara.Range.End

$r = $d.Range($start, $end)
$r.Delete()
